$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.9999978462974438
$ws.Range("D2").Value = 29.99939760023321
$ws.Range("E2").Value = 1.0004
$ws.Range("F2").Value = 30
$ws.Range("C3").Value = 0.9871461082672988
$ws.Range("D3").Value = 29.73463032929338
$ws.Range("E3").Value = 0.9875
$ws.Range("F3").Value = 30.0002
$ws.Range("C4").Value = 0.980243517600477
$ws.Range("D4").Value = 29.5849034281025
$ws.Range("E4").Value = 0.9813
$ws.Range("F4").Value = 30.0077
$ws.Range("C5").Value = 0.9791743678300437
$ws.Range("D5").Value = 29.59300282814636
$ws.Range("E5").Value = 0.98
$ws.Range("F5").Value = 30.0173
$ws.Range("C6").Value = 0.9785802030626806
$ws.Range("D6").Value = 29.57767244903478
$ws.Range("E6").Value = 0.9795
$ws.Range("F6").Value = 30.0005
$ws.Range("C7").Value = 0.9744400820520071
$ws.Range("D7").Value = 29.45855205855314
$ws.Range("E7").Value = 0.9756
$ws.Range("F7").Value = 30.0172
$ws.Range("C8").Value = 0.9735033024337679
$ws.Range("D8").Value = 29.46687977199923
$ws.Range("E8").Value = 0.9728
$ws.Range("F8").Value = 29.986
$ws.Range("C9").Value = 0.9699838927662394
$ws.Range("D9").Value = 29.43931490686414
$ws.Range("E9").Value = 0.9693000000000001
$ws.Range("F9").Value = 29.9969
$ws.Range("C10").Value = 0.9765629813475247
$ws.Range("D10").Value = 29.5689951013718
$ws.Range("E10").Value = 0.9775
$ws.Range("F10").Value = 30.0007
$ws.Range("C11").Value = 0.9764177383884384
$ws.Range("D11").Value = 29.57007168200554
$ws.Range("E11").Value = 0.9772999999999999
$ws.Range("F11").Value = 30.0001
$ws.Range("C12").Value = 0.9742835493041643
$ws.Range("D12").Value = 29.5596200457214
$ws.Range("E12").Value = 0.9752999999999999
$ws.Range("F12").Value = 30.0043
$ws.Range("C13").Value = 0.9787364070787736
$ws.Range("D13").Value = 29.59636570631315
$ws.Range("E13").Value = 0.9772999999999999
$ws.Range("F13").Value = 29.958
$ws.Range("C14").Value = 0.9788612320651877
$ws.Range("D14").Value = 29.59537015192945
$ws.Range("E14").Value = 0.9778
$ws.Range("F14").Value = 29.9668
$ws.Range("C15").Value = 0.9740601048147353
$ws.Range("D15").Value = 29.55886036234379
$ws.Range("E15").Value = 0.975
$ws.Range("F15").Value = 30.0026
$ws.Range("C16").Value = 0.9738772004839216
$ws.Range("D16").Value = 29.56050756558678
$ws.Range("E16").Value = 0.9743000000000001
$ws.Range("F16").Value = 29.9926
$ws.Range("C17").Value = 0.9714830144298983
$ws.Range("D17").Value = 29.58282653418424
$ws.Range("E17").Value = 0.9712
$ws.Range("F17").Value = 30.003
$ws.Range("C18").Value = 0.9710137207108601
$ws.Range("D18").Value = 29.58702900607485
$ws.Range("E18").Value = 0.9696
$ws.Range("F18").Value = 29.984
$ws.Range("C19").Value = 0.9711917221779645
$ws.Range("D19").Value = 29.58522093023345
$ws.Range("E19").Value = 0.9705
$ws.Range("F19").Value = 29.9953
$ws.Range("C20").Value = 0.9667122234852344
$ws.Range("D20").Value = 29.42532072568056
$ws.Range("E20").Value = 0.9659
$ws.Range("F20").Value = 30.0008
$ws.Range("C21").Value = 0.9649053206510023
$ws.Range("D21").Value = 29.4174950943476
$ws.Range("E21").Value = 0.9641
$ws.Range("F21").Value = 29.9995
$ws.Range("C22").Value = 0.9665016089778371
$ws.Range("D22").Value = 29.42723082241818
$ws.Range("E22").Value = 0.9655
$ws.Range("F22").Value = 29.9972
$ws.Range("C23").Value = 0.9686392031826067
$ws.Range("D23").Value = 29.43352862348862
$ws.Range("E23").Value = 0.9678
$ws.Range("F23").Value = 30.0002
$ws.Range("C24").Value = 0.9682178474089123
$ws.Range("D24").Value = 29.43154632128718
$ws.Range("E24").Value = 0.9668
$ws.Range("F24").Value = 29.9993
$ws.Range("C25").Value = 0.9611383124065676
$ws.Range("D25").Value = 29.41404149714652
$ws.Range("E25").Value = 0.9595
$ws.Range("F25").Value = 30.0083
$ws.Range("C26").Value = 0.9608724950534766
$ws.Range("D26").Value = 29.41608971019862
$ws.Range("E26").Value = 0.9583
$ws.Range("F26").Value = 29.9866
$ws.Range("C27").Value = 0.9603569472568562
$ws.Range("D27").Value = 29.4210515907283
$ws.Range("E27").Value = 0.9553
$ws.Range("F27").Value = 29.9427
$ws.Range("C28").Value = 0.9585940595682791
$ws.Range("D28").Value = 29.39090984292467
$ws.Range("E28").Value = 0.9586
$ws.Range("F28").Value = 30.0007
$ws.Range("C29").Value = 0.958450173813186
$ws.Range("D29").Value = 29.39028470881522
$ws.Range("E29").Value = 0.9585
$ws.Range("F29").Value = 30.0012
$ws.Range("C30").Value = 0.9583275399564768
$ws.Range("D30").Value = 29.39296893393077
$ws.Range("E30").Value = 0.9581
$ws.Range("F30").Value = 29.9973
$ws.Range("C31").Value = 0.9757697657743152
$ws.Range("D31").Value = 29.57497001850548
$ws.Range("E31").Value = 0.9761
$ws.Range("F31").Value = 29.9908
$ws.Range("C32").Value = 0.9729840497728343
$ws.Range("D32").Value = 29.46512336021737
$ws.Range("E32").Value = 0.9727
$ws.Range("F32").Value = 30.0104
$ws.Range("C33").Value = 0.9621463298905579
$ws.Range("D33").Value = 29.40598275523224
$ws.Range("E33").Value = 0.9615
$ws.Range("F33").Value = 30.0033
$ws.Range("C34").Value = 0.9599397716136483
$ws.Range("D34").Value = 29.39672884681385
$ws.Range("E34").Value = 0.9597
$ws.Range("F34").Value = 30.0003
$ws.Range("C35").Value = 0.9590405345349261
$ws.Range("D35").Value = 29.3929854208738
$ws.Range("E35").Value = 0.9589
$ws.Range("F35").Value = 30.0003
$ws.Range("C36").Value = 0.9725915297809088
$ws.Range("D36").Value = 29.46867422021144
$ws.Range("E36").Value = 0.9715
$ws.Range("F36").Value = 29.9897
$ws.Range("C37").Value = 0.9728065918307963
$ws.Range("D37").Value = 29.46660637999487
$ws.Range("E37").Value = 0.9708
$ws.Range("F37").Value = 29.9698
